$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows whose only change is "Taxonsorteringsordning" (column B) += 1 ---
$ws.Cells.Item(2, 2).Value = 79244
$ws.Cells.Item(3, 2).Value = 80309
$ws.Cells.Item(4, 2).Value = 83224
$ws.Cells.Item(5, 2).Value = 78256
$ws.Cells.Item(6, 2).Value = 83224
$ws.Cells.Item(8, 2).Value = 79244
$ws.Cells.Item(9, 2).Value = 83224
$ws.Cells.Item(10, 2).Value = 79244
$ws.Cells.Item(12, 2).Value = 91772
$ws.Cells.Item(13, 2).Value = 78256
$ws.Cells.Item(14, 2).Value = 83224
$ws.Cells.Item(15, 2).Value = 83224
$ws.Cells.Item(16, 2).Value = 78256
$ws.Cells.Item(19, 2).Value = 78256
$ws.Cells.Item(20, 2).Value = 83224
$ws.Cells.Item(21, 2).Value = 80315
$ws.Cells.Item(22, 2).Value = 91772
$ws.Cells.Item(23, 2).Value = 83224
$ws.Cells.Item(24, 2).Value = 78256
$ws.Cells.Item(25, 2).Value = 79244
$ws.Cells.Item(26, 2).Value = 91772
$ws.Cells.Item(27, 2).Value = 83224
$ws.Cells.Item(28, 2).Value = 79244
$ws.Cells.Item(29, 2).Value = 83224
$ws.Cells.Item(30, 2).Value = 79244

# --- Rows 17 and 18: the two observation records swapped places ---
# (with column B keeping the "+1" pattern relative to the record now on that row)
$ws.Cells.Item(17, 1).Value = 130864519
$ws.Cells.Item(17, 2).Value = 83224
$ws.Cells.Item(17, 5).Value = 6440
$ws.Cells.Item(17, 6).Value = "Vitgrynig nållav"
$ws.Cells.Item(17, 7).Value = "Chaenotheca subroscida"
$ws.Cells.Item(17, 8).Value = "(Eitner) Zahlbr."
$ws.Cells.Item(17, 17).Value = 446074
$ws.Cells.Item(17, 18).Value = 7030848
$ws.Cells.Item(17, 26).Value = "13:06"
$ws.Cells.Item(17, 28).Value = "13:06"

$ws.Cells.Item(18, 1).Value = 130864522
$ws.Cells.Item(18, 2).Value = 78256
$ws.Cells.Item(18, 5).Value = 228579
$ws.Cells.Item(18, 6).Value = "Liten svartspik"
$ws.Cells.Item(18, 7).Value = "Chaenothecopsis nana"
$ws.Cells.Item(18, 8).Value = "Tibell"
$ws.Cells.Item(18, 17).Value = 446025
$ws.Cells.Item(18, 18).Value = 7031011
$ws.Cells.Item(18, 26).Value = "13:42"
$ws.Cells.Item(18, 28).Value = "13:42"
